$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Monday hours for week of 2018-02-19 (row 6) from 8.25 to 9.25
$ws.Range("B6").Value = 9.25

# Add a grand total row summing the weekly totals
$ws.Range("I19").Formula = "=SUM(I2:I18)"

# Move the active selection to reflect where the user ended up after entering the total
$ws.Range("I20").Select()
